$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add emoji to the section heading "Course Information"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Course Information", $true, $false, $false, $false, $false, $true, 1, $false, "🏛️ Course Information", 1) | Out-Null

# ---------------------------------------------------------------------------
# 2. Remove the whole "Zoom Hours: Friday, 1-3 (schedule online)" paragraph
#    (heading run + body text + hyperlink), which sits between "Office Hours"
#    and "Phone" in the Instructor block.
# ---------------------------------------------------------------------------
$rng = $d.Content
if ($rng.Find.Execute("Zoom Hours")) {
    $zoomPara = $rng.Paragraphs(1)
    $zoomPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3. Add a new "Quick Description" section right after the Instructor
#    section (after the "e-mail" paragraph), before "Course Description".
#    Insert the break+text from a point inside the e-mail paragraph (not at
#    its exact end) so the new paragraphs land before the existing
#    "instructor" bookmark's end marker rather than after it.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("e-mail") | Out-Null
$emailPara = $rng2.Paragraphs(1)
$endPos = $emailPara.Range.End
$insertRng = $d.Range($endPos - 1, $endPos - 1)
$insertRng.InsertAfter("`rQuick Description`rThis course is designed for students who want to think critically about artificial intelligence and design artificial intelligence policies for government, business, and nonprofit organizations of all sizes.")

# Re-fetch the two freshly-created paragraphs (the two paragraphs right after
# the now-unchanged e-mail paragraph) and give them the right styles.
$afterEmailRng = $d.Range($emailPara.Range.End, $d.Content.End)
$headingPara = $afterEmailRng.Paragraphs(1)
$bodyPara = $afterEmailRng.Paragraphs(2)
$headingPara.Style = "Heading3"
$bodyPara.Style = "FirstParagraph"

# Wrap the two new paragraphs in a "quick-description" bookmark.
$bmRange = $d.Range($headingPara.Range.Start, $bodyPara.Range.End)
$d.Bookmarks.Add("quick-description", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 4. Update section headings with their emoji prefixes.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Materials", $true, $false, $false, $false, $false, $true, 1, $false, "📘 Materials", 1) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Course Requirements", $true, $false, $false, $false, $false, $true, 1, $false, "📌 Course Requirements", 1) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Course Expectations & Guidelines", $true, $false, $false, $false, $false, $true, 1, $false, "🎓 Course Expectations & Guidelines", 1) | Out-Null

$rng = $d.Content
$rng.Find.Execute("Learning Objective Assessment", $true, $false, $false, $false, $false, $true, 1, $false, "📊 Learning Objective Assessment", 1) | Out-Null

# ---------------------------------------------------------------------------
# 5. Tweak the "E-mail" etiquette paragraph wording.
# ---------------------------------------------------------------------------
$rng = $d.Content
$oldEmailPara = "Email is the best way to contact me. I’m usually pretty responsive, but as a baseline, I always aim to get back to you in a modified 24-hour fashion: by the end of the business day the day after you email. So if you email me at 2 PM Tuesday, I’ll get back to you by 6 PM Wednesday at the latest; if 10 PM Thursday, by 6 PM Friday; if you email me at 3 PM on Friday, by 6 PM Monday, etc."
$newEmailPara = "Email is the best way to contact me. I try to be pretty responsive, but as a baseline, I always aim to get back to you in a modified 24-hour fashion: by the end of the business day the day after you email (at minimum). So if you email me at 2 PM Tuesday, I’ll get back to you by 6 PM Wednesday at the latest; if 10 PM Thursday, by 6 PM Friday; if you email me at 3 PM on Friday, by 6 PM Monday, etc."
$rng.Find.Execute($oldEmailPara, $true, $false, $false, $false, $false, $true, 1, $false, $newEmailPara, 1) | Out-Null

# ---------------------------------------------------------------------------
# 6. Extend the footnote about response-time reminders.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Footnotes.Count; $i++) {
    $fn = $d.Footnotes.Item($i)
    if ($fn.Range.Text -like "Usually I aim to be much much faster*") {
        $fn.Range.Text = $fn.Range.Text.Replace(
            "feel free to bump a reminder.",
            "feel free to bump a reminder. No hard feelings. Sometimes things get busy and I lose track of an email."
        )
        break
    }
}
